$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 34 (the empty cell below the
# "Abläufe..." header), shifting everything below it down by one.
$ws.Rows.Item(34).Insert()

# New explanatory bullet point, styled like the bold header in C33
# (leading apostrophe forces the text quote-prefix, matching C33's style).
$ws.Range("C34").Value = "'-> wenn es auch mit Webclient funktionieren muss, müssen die genutzten öffentlichen Prozeduren kopiert werden als custom version und aus diesen die Prüfung des Transactusername und die Prüfung der Schreibrechte entfernt werden"
$ws.Range("C34").Font.Bold = $true

# Move the view / selection roughly where the author left it.
$ws.Application.ActiveWindow.ScrollRow = 10
$ws.Range("C35").Select() | Out-Null
